# The document embeds the same two logos (Pearson "PearsonLogo.png" and
# "BTec_Logo-Orange") in its headers/footers, but their drawing objects'
# display names were swapped relative to the actual media part filenames:
#   - footer1.xml / footer2.xml: the Pearson logo's drawing is named
#     "image1.png" and should be renamed "image2.png"
#   - header1.xml: the BTec logo's drawing is named "image2.jpg" and
#     should be renamed "image1.jpg"
# Both the <wp:docPr> and the inner <pic:cNvPr> "name" attribute need to
# be updated for each picture. InlineShape.Name only round-trips through
# <wp:docPr>, so instead we patch the document's raw OOXML text, which
# keeps every other byte of the package untouched.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
